$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (unchanged text, kept here for completeness)
$ws.Range("A1").Value = "SKU"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Unit Cost"
$ws.Range("E1").Value = "Pack. Cost"
$ws.Range("F1").Value = "Ship. Cost"
$ws.Range("G1").Value = "Market Fee"
$ws.Range("H1").Value = "VAT(%)"

# Row 2 (existing row, category/name text stays literal string values)
$ws.Range("A2").Value = "SKU101"
$ws.Range("B2").Value = "Food Container"
$ws.Range("C2").Value = "Smile Container"
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 2.5
$ws.Range("F2").Value = 1.5
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 5

# New rows 3-10
$ws.Range("A3").Value = "SKU102"
$ws.Range("B3").Value = "Water Boottle"
$ws.Range("C3").Value = "PolyProphelene"
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1.5
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3

$ws.Range("A4").Value = "SKU103"
$ws.Range("B4").Value = "TableWare"
$ws.Range("C4").Value = "Vigoset"
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 2.5
$ws.Range("F4").Value = 1.5
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 3.5

$ws.Range("A5").Value = "SKU104"
$ws.Range("B5").Value = "KitchenWare"
$ws.Range("C5").Value = "DishRack"
$ws.Range("D5").Value = 34
$ws.Range("E5").Value = 23
$ws.Range("F5").Value = 44
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 2

$ws.Range("A6").Value = "SKU105"
$ws.Range("B6").Value = "Hanger"
$ws.Range("C6").Value = "Hanger"
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = 0.25
$ws.Range("F6").Value = 0.15
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0.75

$ws.Range("A7").Value = "SKU106"
$ws.Range("B7").Value = "Cleaning"
$ws.Range("C7").Value = "Flip & Clean"
$ws.Range("D7").Value = 60
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 7.5

$ws.Range("A8").Value = "SKU107"
$ws.Range("B8").Value = "Storage"
$ws.Range("C8").Value = "Rack"
$ws.Range("D8").Value = 250
$ws.Range("E8").Value = 25
$ws.Range("F8").Value = 35
$ws.Range("G8").Value = 20
$ws.Range("H8").Value = 9.3

$ws.Range("A9").Value = "SKU108"
$ws.Range("B9").Value = "Furniture"
$ws.Range("C9").Value = "Table"
$ws.Range("D9").Value = 175
$ws.Range("E9").Value = 30
$ws.Range("F9").Value = 18.75
$ws.Range("G9").Value = 15
$ws.Range("H9").Value = 7

$ws.Range("A10").Value = "SKU109"
$ws.Range("B10").Value = "Kids"
$ws.Range("C10").Value = "Tiffin Box"
$ws.Range("D10").Value = 35
$ws.Range("E10").Value = 1.5
$ws.Range("F10").Value = 0.15
$ws.Range("G10").Value = 0.53
$ws.Range("H10").Value = 0.15

# Match final selection state from the diff
$ws.Range("D10").Select()
